$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.5690543055534363
$ws.Range("B1").Value = 1.994382739067078
$ws.Range("C1").Value = 6.316158771514893
$ws.Range("D1").Value = 2.715606927871704
$ws.Range("E1").Value = 1.851594924926758
